# Auto-generated edit script
# Applies updated game-simulation probability matrix values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2087542087542087
$ws.Range("C2").Value = 0.531986531986532
$ws.Range("J2").Value = 0.02356902356902357
$ws.Range("O2").Value = 0.003367003367003367
$ws.Range("P2").Value = 0.1346801346801347
$ws.Range("S2").Value = 0.09764309764309764
$ws.Range("B3").Value = 0.01851851851851852
$ws.Range("C3").Value = 0.03703703703703703
$ws.Range("J3").Value = 0.01851851851851852
$ws.Range("P3").Value = 0.7283950617283951
$ws.Range("S3").Value = 0.1975308641975309
$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.6808510638297872
$ws.Range("S4").Value = 0.2978723404255319
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.04583333333333333
$ws.Range("D6").Value = 0.008333333333333333
$ws.Range("E6").Value = 0.004166666666666667
$ws.Range("F6").Value = 0.08333333333333333
$ws.Range("J6").Value = 0.2291666666666667
$ws.Range("O6").Value = 0.02916666666666667
$ws.Range("Q6").Value = 0.1708333333333333
$ws.Range("R6").Value = 0.0375
$ws.Range("S6").Value = 0.3916666666666667
$ws.Range("B7").Value = 0.1052631578947368
$ws.Range("D7").Value = 0.02923976608187134
$ws.Range("F7").Value = 0.1052631578947368
$ws.Range("J7").Value = 0.1169590643274854
$ws.Range("O7").Value = 0.03508771929824561
$ws.Range("Q7").Value = 0.1637426900584795
$ws.Range("R7").Value = 0.07017543859649122
$ws.Range("S7").Value = 0.3742690058479532
$ws.Range("B8").Value = 0.08884297520661157
$ws.Range("D8").Value = 0.02272727272727273
$ws.Range("F8").Value = 0.0640495867768595
$ws.Range("J8").Value = 0.09710743801652892
$ws.Range("O8").Value = 0.03099173553719008
$ws.Range("Q8").Value = 0.1528925619834711
$ws.Range("R8").Value = 0.08264462809917356
$ws.Range("S8").Value = 0.4607438016528926
$ws.Range("B9").Value = 0.09917355371900827
$ws.Range("D9").Value = 0.02066115702479339
$ws.Range("F9").Value = 0.07851239669421488
$ws.Range("J9").Value = 0.09917355371900827
$ws.Range("O9").Value = 0.04958677685950413
$ws.Range("Q9").Value = 0.1487603305785124
$ws.Range("R9").Value = 0.07024793388429752
$ws.Range("S9").Value = 0.4338842975206612
$ws.Range("B10").Value = 0.1044657097288676
$ws.Range("D10").Value = 0.0215311004784689
$ws.Range("E10").Value = 0.0007974481658692185
$ws.Range("F10").Value = 0.07336523125996811
$ws.Range("J10").Value = 0.1172248803827751
$ws.Range("O10").Value = 0.01116427432216906
$ws.Range("Q10").Value = 0.2169059011164274
$ws.Range("R10").Value = 0.07177033492822966
$ws.Range("S10").Value = 0.3827751196172249
$ws.Range("G11").Value = 0.1095890410958904
$ws.Range("J11").Value = 0.1061643835616438
$ws.Range("K11").Value = 0.1883561643835616
$ws.Range("L11").Value = 0.5616438356164384
$ws.Range("S11").Value = 0.03424657534246575
$ws.Range("G12").Value = 0.7062146892655368
$ws.Range("J12").Value = 0.1581920903954802
$ws.Range("K12").Value = 0.01129943502824859
$ws.Range("L12").Value = 0.06779661016949153
$ws.Range("S12").Value = 0.05649717514124294
$ws.Range("G13").Value = 0.5405405405405406
$ws.Range("J13").Value = 0.3783783783783784
$ws.Range("S13").Value = 0.08108108108108109
$ws.Range("F15").Value = 0.01376146788990826
$ws.Range("H15").Value = 0.1651376146788991
$ws.Range("I15").Value = 0.06422018348623854
$ws.Range("J15").Value = 0.3211009174311927
$ws.Range("K15").Value = 0.06880733944954129
$ws.Range("M15").Value = 0.02752293577981652
$ws.Range("N15").Value = 0.004587155963302753
$ws.Range("O15").Value = 0.07798165137614679
$ws.Range("S15").Value = 0.2568807339449541
$ws.Range("F16").Value = 0.02197802197802198
$ws.Range("H16").Value = 0.2417582417582418
$ws.Range("I16").Value = 0.09340659340659341
$ws.Range("J16").Value = 0.4175824175824176
$ws.Range("K16").Value = 0.0989010989010989
$ws.Range("O16").Value = 0.03846153846153846
$ws.Range("S16").Value = 0.08791208791208792
$ws.Range("F17").Value = 0.01552106430155211
$ws.Range("H17").Value = 0.2084257206208426
$ws.Range("I17").Value = 0.1019955654101996
$ws.Range("J17").Value = 0.3991130820399113
$ws.Range("K17").Value = 0.1130820399113082
$ws.Range("M17").Value = 0.01330376940133038
$ws.Range("O17").Value = 0.07538802660753881
$ws.Range("S17").Value = 0.07317073170731707
$ws.Range("F18").Value = 0.01204819277108434
$ws.Range("H18").Value = 0.144578313253012
$ws.Range("I18").Value = 0.1144578313253012
$ws.Range("J18").Value = 0.4397590361445783
$ws.Range("K18").Value = 0.1204819277108434
$ws.Range("M18").Value = 0.02409638554216868
$ws.Range("O18").Value = 0.06024096385542169
$ws.Range("S18").Value = 0.08433734939759036
$ws.Range("F19").Value = 0.01319875776397516
$ws.Range("H19").Value = 0.2243788819875776
$ws.Range("I19").Value = 0.1141304347826087
$ws.Range("J19").Value = 0.3835403726708074
$ws.Range("K19").Value = 0.09937888198757763
$ws.Range("M19").Value = 0.01863354037267081
$ws.Range("N19").Value = 0.001552795031055901
$ws.Range("O19").Value = 0.05279503105590062
$ws.Range("S19").Value = 0.09239130434782608
